$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing user row (row 2) with corrected values
$ws.Range("A2").Value = "saulosjss33"
$ws.Range("G2").Value = "(11)921316555"
$ws.Range("H2").Value = "Brazil"
$ws.Range("L2").Value = "(11)12345"
$ws.Range("M2").Value = "Pass"

# Add new registered user in row 3
$ws.Range("A3").Value = "rapaz"
$ws.Range("B3").Value = "saulojosilva@hotmail.com"
$ws.Range("C3").Value = "mano"
$ws.Range("D3").Value = "eita"
$ws.Range("E3").Value = "Marcos"
$ws.Range("F3").Value = "Silva"
$ws.Range("G3").Value = "(11)954423458"
$ws.Range("H3").Value = "Algeria"
$ws.Range("I3").Value = "paris"
$ws.Range("J3").Value = "rua mundi"
$ws.Range("K3").Value = "af"
$ws.Range("L3").Value = "(11)12345"

# Hyperlink for new user's email, matching style of B2
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:saulojosilva@hotmail.com")
$ws.Range("B3").Style = "Hiperlink"

# Move the active selection to A2 (print/screenshot of final screen)
$ws.Range("A2").Select()
